$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new row of data (Merge Two Sorted Lists)
$ws.Range("A4").Value = "Merge Two Sorted Lists"
$ws.Range("B4").Value = "Easy"
$ws.Range("C4").Value = "Using linkedlists concept, compare node values."

# Match formatting of the rows above (center alignment for A/B, vertical-center for C)
$ws.Range("A4:B4").HorizontalAlignment = -4108
$ws.Range("A4:C4").VerticalAlignment = -4108

# Row height for the new data row
$ws.Rows.Item(4).RowHeight = 28.8

# Update the active selection
$ws.Range("F3").Select()
